$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "24×83=1992" "41×88=3608"
Replace-Text "93×81=7533" "35×80=2800"
Replace-Text "40×42=1680" "75×73=5475"
Replace-Text "92×48=4416" "63×17=1071"
Replace-Text "97×11=1067" "93×99=9207"
Replace-Text "49×73=3577" "25×97=2425"
Replace-Text "36×74=2664" "20×91=1820"
Replace-Text "76×67=5092" "36×91=3276"
Replace-Text "90×85=7650" "32×16=512"
Replace-Text "28×64=1792" "62×23=1426"
Replace-Text "69×94=6486" "15×37=555"
Replace-Text "61×47=2867" "48×32=1536"
Replace-Text "26×89=2314" "88×63=5544"
Replace-Text "33×31=1023" "99×79=7821"
Replace-Text "27×97=2619" "33×92=3036"
Replace-Text "95×46=4370" "25×77=1925"
Replace-Text "54×14=756" "88×74=6512"
Replace-Text "56×20=1120" "44×78=3432"
Replace-Text "48×84=4032" "95×36=3420"
Replace-Text "15×88=1320" "44×13=572"
Replace-Text "69×37=2553" "93×13=1209"
Replace-Text "66×87=5742" "29×85=2465"
Replace-Text "23×38=874" "42×18=756"
Replace-Text "38×45=1710" "16×18=288"
Replace-Text "28×99=2772" "25×35=875"

Write-Output "Done replacing all cells"
